$d = $word.ActiveDocument

$pairs = @(
    @("54×24=", "30×16="),
    @("90×84=", "48×40="),
    @("97×22=", "82×23="),
    @("99×41=", "57×37="),
    @("62×75=", "32×27="),
    @("49×91=", "46×34="),
    @("93×20=", "27×99="),
    @("23×75=", "25×22="),
    @("95×45=", "79×49="),
    @("56×33=", "41×91="),
    @("50×89=", "32×25="),
    @("25×47=", "64×80="),
    @("62×21=", "45×18="),
    @("89×89=", "72×92="),
    @("95×99=", "92×77="),
    @("31×51=", "59×99="),
    @("15×66=", "43×26="),
    @("25×21=", "52×23="),
    @("50×33=", "92×77="),
    @("31×47=", "95×97="),
    @("45×41=", "91×13="),
    @("75×55=", "69×74="),
    @("31×36=", "27×35="),
    @("20×86=", "81×35="),
    @("24×23=", "62×84=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
